# Regenerate the "K" column (strikeouts, column G) values on the active sheet.
# The previous values were based on "Strike#"; the new values are based on "K".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 0
    7  = 2
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 1
    14 = 0
    15 = 2
    16 = 0
    17 = 4
    18 = 2
    19 = 0
    20 = 2
    21 = 2
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 3
    27 = 4
    28 = 2
    29 = 2
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
